# Apply updated cryptos data (values only; keep default/General cell style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.614.49"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -3.08%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.850.91"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -3.70%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -1.15%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'335.80"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +2.96%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'1.002"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.95%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4663"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -3.01%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.3903"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -3.45%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'46.17"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -2.69%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.07913"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -3.36%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.9796"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -2.74%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'22.31"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -6.21%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'1.883.30"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -2.27%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'5.826"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -4.22%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'6.990"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -4.31%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'0.06923"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.71%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'1.002"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -1.05%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'87.62"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -4.26%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  -3.46%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'17.08"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -3.06%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  -1.00%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'28.624.89"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -3.05%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'5.395"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -4.65%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  -5.97%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'2.160"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -0.77%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('B26').Value = "'Monero"
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').Value = "'153.09"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -1.79%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('B27').Value = "'EthereumClassic"
$ws.Range('B27').Style = 'Normal'
$ws.Range('C27').Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range('C27').Style = 'Normal'
$ws.Range('D27').Value = "'19.44"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -2.82%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('B28').Value = "'InternetComputer(DFINITY)"
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range('C28').Style = 'Normal'
$ws.Range('D28').Value = "'6.078"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -4.82%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('B29').Value = "'LidoDAOToken"
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').Value = "'2.026"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -2.70%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('B30').Value = "'BitcoinCash"
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').Value = "'117.42"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -2.48%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value = "'ImmutableX"
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = "'0.9708"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -4.13%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = "'Stellar"
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = "'0.09337"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -2.54%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = "'Filecoin"
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = "'5.361"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -4.21%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value = "'HuobiToken"
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = "'3.482"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -2.22%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('B35').Value = "'ARBITRUM"
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').Value = "'1.345"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -2.84%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('B36').Value = "'Hedera"
$ws.Range('B36').Style = 'Normal'
$ws.Range('C36').Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range('C36').Style = 'Normal'
$ws.Range('D36').Value = "'0.06168"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -2.89%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('B37').Value = "'VeChain"
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').Value = "'0.02202"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -3.59%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('B38').Value = "'TrustWalletToken"
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = "'1.170"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.76%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = "'TheSandbox"
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = "'0.5711"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -3.84%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = "'FraxShare"
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = "'7.683"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -2.57%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = "'Aptos"
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = "'10.15"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -5.34%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = "'Algorand"
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = "'0.1791"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -2.72%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = "'RenderToken"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'2.409"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -2.89%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'WEMIXToken"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'1.220"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -4.39%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = "'Decentraland"
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = "'0.5373"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -3.04%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = "'EnergySwap"
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'11.75"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -5.30%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = "'Cronos"
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = "'0.07100"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -4.93%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = "'NEARProtocol"
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = "'1.905"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -3.42%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = "'Quant"
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'113.39"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -3.71%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = "'PaxDollar"
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'1.001"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -1.02%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'MXToken"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'2.336"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -3.90%  "
$ws.Range('E51').Style = 'Normal'
